$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($addr, $val)
    $ws.Range($addr).Value = "'" + $val
    $ws.Range($addr).Style = "Normal"
}

# Row 2
Set-TextCell "D2" "28.243.05"
Set-TextCell "E2" "  +0.95%  "

# Row 3
Set-TextCell "D3" "1.884.10"
Set-TextCell "E3" "  +1.09%  "

# Row 4
Set-TextCell "E4" "  +0.09%  "

# Row 5
Set-TextCell "D5" "314.11"
Set-TextCell "E5" "  +0.64%  "

# Row 6
Set-TextCell "E6" "  +0.13%  "

# Row 7
Set-TextCell "D7" "0.5148"
Set-TextCell "E7" "  +1.36%  "

# Row 8
Set-TextCell "D8" "0.3903"
Set-TextCell "E8" "  +2.90%  "

# Row 9
Set-TextCell "D9" "0.08370"
Set-TextCell "E9" "  +0.71%  "

# Row 10
Set-TextCell "D10" "1.121"
Set-TextCell "E10" "  +1.45%  "

# Row 11
Set-TextCell "E11" "  +0.48%  "

# Row 12
Set-TextCell "D12" "6.246"
Set-TextCell "E12" "  +0.84%  "

# Row 13
Set-TextCell "B13" "WrappedEther"
Set-TextCell "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextCell "D13" "1.897.20"
Set-TextCell "E13" "  +1.72%  "

# Row 14
Set-TextCell "B14" "Solana"
Set-TextCell "C14" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextCell "D14" "20.75"
Set-TextCell "E14" "  +1.60%  "

# Row 15
Set-TextCell "D15" "7.296"
Set-TextCell "E15" "  +1.80%  "

# Row 16
Set-TextCell "D16" "1.006"
Set-TextCell "E16" "  +0.08%  "

# Row 17
Set-TextCell "D17" "0.00001108"
Set-TextCell "E17" "  +1.39%  "

# Row 18
Set-TextCell "D18" "91.46"
Set-TextCell "E18" "  +1.32%  "

# Row 19
Set-TextCell "D19" "0.06667"
Set-TextCell "E19" "  +0.77%  "

# Row 20
Set-TextCell "D20" "17.85"
Set-TextCell "E20" "  +0.23%  "

# Row 21
Set-TextCell "E21" "  +0.24%  "

# Row 22
Set-TextCell "D22" "6.078"
Set-TextCell "E22" "  +1.21%  "

# Row 23
Set-TextCell "D23" "28.279.82"
Set-TextCell "E23" "  +0.93%  "

# Row 24
Set-TextCell "D24" "11.18"
Set-TextCell "E24" "  +0.84%  "

# Row 25
Set-TextCell "E25" "  +0.66%  "

# Row 26
Set-TextCell "D26" "2.097.58"
Set-TextCell "E26" "  +0.79%  "

# Row 27
Set-TextCell "D27" "2.520"
Set-TextCell "E27" "  -1.76%  "

# Row 28
Set-TextCell "D28" "159.10"
Set-TextCell "E28" "  +1.17%  "

# Row 29
Set-TextCell "D29" "20.67"
Set-TextCell "E29" "  +1.46%  "

# Row 30
Set-TextCell "D30" "125.61"
Set-TextCell "E30" "  -0.50%  "

# Row 31
Set-TextCell "D31" "0.1067"
Set-TextCell "E31" "  +1.28%  "

# Row 32
Set-TextCell "D32" "1.047"
Set-TextCell "E32" "  +1.12%  "

# Row 33
Set-TextCell "D33" "5.890"
Set-TextCell "E33" "  +5.56%  "

# Row 34
Set-TextCell "D34" "3.602"
Set-TextCell "E34" "  +0.21%  "

# Row 35
Set-TextCell "D35" "9.786"
Set-TextCell "E35" "  +1.64%  "

# Row 36
Set-TextCell "E36" "  +1.77%  "

# Row 37
Set-TextCell "D37" "0.06582"
Set-TextCell "E37" "  +1.14%  "

# Row 38
Set-TextCell "E38" "  +2.08%  "

# Row 39
Set-TextCell "E39" "  +0.76%  "

# Row 40
Set-TextCell "D40" "0.6547"
Set-TextCell "E40" "  +2.86%  "

# Row 41
Set-TextCell "D41" "5.026"
Set-TextCell "E41" "  +3.60%  "

# Row 42
Set-TextCell "E42" "  +0.14%  "

# Row 43
Set-TextCell "E43" "  +0.73%  "

# Row 44
Set-TextCell "D44" "0.6139"
Set-TextCell "E44" "  +1.24%  "

# Row 45
Set-TextCell "D45" "13.12"
Set-TextCell "E45" "  +1.25%  "

# Row 46
Set-TextCell "D46" "1.289"
Set-TextCell "E46" "  +0.38%  "

# Row 47
Set-TextCell "D47" "3.680"
Set-TextCell "E47" "  +0.56%  "

# Row 48
Set-TextCell "B48" "EOS"
Set-TextCell "C48" "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextCell "D48" "1.239"
Set-TextCell "E48" "  +2.50%  "

# Row 49
Set-TextCell "B49" "NEARProtocol"
Set-TextCell "C49" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell "D49" "2.014"
Set-TextCell "E49" "  +1.47%  "

# Row 50
Set-TextCell "D50" "121.57"
Set-TextCell "E50" "  +0.44%  "

# Row 51
Set-TextCell "D51" "79.08"
Set-TextCell "E51" "  -0.70%  "
